$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "pass1"
$ws.Range("C3").Value = "pass2"
$ws.Range("C4").Value = "pass3"
$ws.Range("B5").Value = "juani"
$ws.Range("C5").Value = "pass4"
